$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $new, 2) | Out-Null
}

Replace-Text "2023-08-10 Thursday" "2023-08-11 Friday"

Replace-Text "88×42=3696" "28×31=868"
Replace-Text "13×39=507" "52×43=2236"
Replace-Text "57×28=1596" "12×89=1068"
Replace-Text "52×58=3016" "48×41=1968"
Replace-Text "43×90=3870" "15×70=1050"

Replace-Text "26×73=1898" "54×95=5130"
Replace-Text "62×59=3658" "39×46=1794"
Replace-Text "97×71=6887" "85×36=3060"
Replace-Text "40×91=3640" "83×89=7387"
Replace-Text "38×45=1710" "13×27=351"

Replace-Text "98×62=6076" "86×42=3612"
Replace-Text "77×45=3465" "63×76=4788"
Replace-Text "84×86=7224" "43×82=3526"
Replace-Text "42×28=1176" "81×89=7209"
Replace-Text "47×29=1363" "40×54=2160"

Replace-Text "80×47=3760" "14×44=616"
Replace-Text "29×37=1073" "30×40=1200"
Replace-Text "65×58=3770" "15×95=1425"
Replace-Text "99×80=7920" "72×97=6984"
Replace-Text "43×86=3698" "39×57=2223"

Replace-Text "85×28=2380" "42×52=2184"
Replace-Text "25×91=2275" "38×88=3344"
Replace-Text "91×32=2912" "21×69=1449"
Replace-Text "77×53=4081" "11×38=418"
Replace-Text "19×75=1425" "39×89=3471"
